# upd yeast growth data
# - add a second relative-growth-rate column (RGR2, log base 2 variant)
# - rename the DL*/LL* condition labels to DLA*/LLA* (keep the same meaning)
# - append a new sample (20230808_0: control / DLA4500 / DLS4500)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column I ---
$ws.Range("I1").Value = "RGR2"

# --- Rename condition labels: DL0005/DL0050/DL0500/DL4500/LL4500 -> DLA.../LLA... ---
$renameMap = @{
    "DL0005" = "DLA0005"
    "DL0050" = "DLA0050"
    "DL0500" = "DLA0500"
    "DL4500" = "DLA4500"
    "LL4500" = "LLA4500"
}
for ($r = 2; $r -le 43; $r++) {
    $cur = $ws.Cells.Item($r, 2).Value()
    if ($renameMap.ContainsKey($cur)) {
        $ws.Cells.Item($r, 2).Value = $renameMap[$cur]
    }
}

# --- New RGR2 column (log2-based relative growth rate) for every existing data row ---
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=(LOG(D$r,2)-LOG(C$r,2))/E$r"
}

# --- New sample rows 44-46 (20230808_0) ---
$ws.Range("A44").Value = "20230808_0"
$ws.Range("B44").Value = "control"
$ws.Range("C44").Value = 0.433
$ws.Range("D44").Value = 1.037
$ws.Range("E44").Formula = "=3+35/60"
$ws.Range("F44").Formula = "=LOG(D44/C44)"
$ws.Range("G44").Formula = "=E44*LOG(2)/F44"
$ws.Range("H44").Formula = "=(LN(D44)-LN(C44))/E44"
$ws.Range("I44").Formula = "=(LOG(D44,2)-LOG(C44,2))/E44"

$ws.Range("A45").Value = "20230808_0"
$ws.Range("B45").Value = "DLA4500"
$ws.Range("C45").Value = 0.436
$ws.Range("D45").Value = 0.838
$ws.Range("E45").Formula = "=3+35/60"
$ws.Range("F45").Formula = "=LOG(D45/C45)"
$ws.Range("G45").Formula = "=E45*LOG(2)/F45"
$ws.Range("H45").Formula = "=(LN(D45)-LN(C45))/E45"
$ws.Range("I45").Formula = "=(LOG(D45,2)-LOG(C45,2))/E45"

$ws.Range("A46").Value = "20230808_0"
$ws.Range("B46").Value = "DLS4500"
$ws.Range("C46").Value = 0.527
$ws.Range("D46").Value = 1.211
$ws.Range("E46").Formula = "=3+35/60"
$ws.Range("F46").Formula = "=LOG(D46/C46)"
$ws.Range("G46").Formula = "=E46*LOG(2)/F46"
$ws.Range("H46").Formula = "=(LN(D46)-LN(C46))/E46"
$ws.Range("I46").Formula = "=(LOG(D46,2)-LOG(C46,2))/E46"

# --- Sheet view: scroll down a bit, select D44 (matches the updated selection in the source) ---
$ws.Range("A31").Select()
$ws.Range("D44").Select()
